$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.471.44"
$ws.Range("D3").Value = "2.988.13"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.74"
$ws.Range("E5").Value = "  +3.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.32"
$ws.Range("E6").Value = "  +2.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  +1.70%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.26"
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("D13").Value = "3.457.57"
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.42"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").Value = "2.997.38"
$ws.Range("E16").Value = "  +3.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.974"
$ws.Range("E17").Value = "  +6.13%  "
$ws.Range("D18").Value = "51.457.04"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.32"
$ws.Range("E19").Value = "  +4.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.44"
$ws.Range("E20").Value = "  +4.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.94"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("E22").Value = "  +2.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.17"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.96"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.93"
$ws.Range("E25").Value = "  +10.41%  "
$ws.Range("E26").Value = "  +17.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.77"
$ws.Range("E27").Value = "  +25.18%  "
$ws.Range("E28").Value = "  +15.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.170"
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.00"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.87"
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.01"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("E36").Value = "  +8.85%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.11"
$ws.Range("E39").Value = "  +1.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.58"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.84"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  +4.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.41"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.80"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  +19.06%  "
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("E47").Value = "  +2.97%  "
$ws.Range("D48").Value = "2.033.13"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  +5.11%  "
$ws.Range("E50").Value = "  +8.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.24"
$ws.Range("E51").Value = "  +3.94%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D51").Style = "Normal"
